{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values per data-row (0-based row indices 0, 4, 8, 12, 16), 5 columns each.\nconst newValues = [\n  [0, [\"289\u00f76=48, 1\", \"934\u00f72=467, 0\", \"173\u00f78=21, 5\", \"372\u00f79=41, 3\", \"415\u00f74=103, 3\"]],\n  [4, [\"498\u00f75=99, 3\", \"986\u00f73=328, 2\", \"339\u00f76=56, 3\", \"382\u00f78=47, 6\", \"973\u00f78=121, 5\"]],\n  [8, [\"847\u00f73=282, 1\", \"266\u00f79=29, 5\", \"603\u00f76=100, 3\", \"203\u00f74=50, 3\", \"391\u00f77=55, 6\"]],\n  [12, [\"286\u00f74=71, 2\", \"456\u00f76=76, 0\", \"836\u00f75=167, 1\", \"942\u00f75=188, 2\", \"312\u00f78=39, 0\"]],\n  [16, [\"201\u00f74=50, 1\", \"564\u00f74=141, 0\", \"509\u00f72=254, 1\", \"972\u00f78=121, 4\", \"800\u00f74=200, 0\"]],\n];\n\nfor (const [rowIndex, rowValues] of newValues) {\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n    range.insertText(rowValues[col], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"289\u00f76=48, 1\", \"934\u00f72=467, 0\", \"173\u00f78=21, 5\", \"372\u00f79=41, 3\", \"415\u00f74=103, 3\"),\n    @(\"498\u00f75=99, 3\", \"986\u00f73=328, 2\", \"339\u00f76=56, 3\", \"382\u00f78=47, 6\", \"973\u00f78=121, 5\"),\n    @(\"847\u00f73=282, 1\", \"266\u00f79=29, 5\", \"603\u00f76=100, 3\", \"203\u00f74=50, 3\", \"391\u00f77=55, 6\"),\n    @(\"286\u00f74=71, 2\", \"456\u00f76=76, 0\", \"836\u00f75=167, 1\", \"942\u00f75=188, 2\", \"312\u00f78=39, 0\"),\n    @(\"201\u00f74=50, 1\", \"564\u00f74=141, 0\", \"509\u00f72=254, 1\", \"972\u00f78=121, 4\", \"800\u00f74=200, 0\")\n)\n\n$rows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $rows.Length; $i++) {\n    $r = $rows[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i][$c - 1]\n    }\n}\n"}
